$wb = $excel.ActiveWorkbook

# --- Add the new "Colors" worksheet at the end of the tab order ---
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$ws = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$ws.Name = "Colors"

# --- Header row ---
$ws.Cells.Item(1, 1).Value = "key"
$ws.Cells.Item(1, 2).Value = "en"
$ws.Cells.Item(1, 3).Value = "es"
$ws.Cells.Item(1, 4).Value = "it"
$ws.Cells.Item(1, 5).Value = "fr"
$ws.Cells.Item(1, 6).Value = "de"
$ws.Cells.Item(1, 7).Value = "nl"
$ws.Cells.Item(1, 8).Value = "ja"
$ws.Cells.Item(1, 9).Value = "fa"
$ws.Cells.Item(1, 10).Value = "notes"

# --- Data rows: group key, color swatches (hex) per locale, notes ---
# Row 2 - A - pale green
$ws.Cells.Item(2, 1).Value = "A"
$ws.Cells.Item(2, 2).Value = "#b0df90"
$ws.Cells.Item(2, 3).Value = "#97cd73"
$ws.Cells.Item(2, 10).Value = "pale green"

# Row 3 - B - pale cyan
$ws.Cells.Item(3, 1).Value = "B"
$ws.Cells.Item(3, 2).Value = "#92d6df"
$ws.Cells.Item(3, 10).Value = "pale cyan"

# Row 4 - C - pale red
$ws.Cells.Item(4, 1).Value = "C"
$ws.Cells.Item(4, 2).Value = "#f79d8f"
$ws.Cells.Item(4, 3).Value = "#e2694f"
$ws.Cells.Item(4, 4).Value = "#eb9784"
$ws.Cells.Item(4, 10).Value = "pale red"

# Row 5 - D - pale pink
$ws.Cells.Item(5, 1).Value = "D"
$ws.Cells.Item(5, 2).Value = "#efc9d8"
$ws.Cells.Item(5, 10).Value = "pale pink"

# Row 6 - E - pale rose
$ws.Cells.Item(6, 1).Value = "E"
$ws.Cells.Item(6, 2).Value = "#eb84af"
$ws.Cells.Item(6, 3).Value = "#f173a7"
$ws.Cells.Item(6, 10).Value = "pale rose"

# Row 7 - F - pale teal
$ws.Cells.Item(7, 1).Value = "F"
$ws.Cells.Item(7, 2).Value = "#60b4c7"
$ws.Cells.Item(7, 3).Value = "#418e9f"
$ws.Cells.Item(7, 10).Value = "pale teal"

# Row 8 - G - pale orange
$ws.Cells.Item(8, 1).Value = "G"
$ws.Cells.Item(8, 2).Value = "#fab077"
$ws.Cells.Item(8, 3).Value = "#faae53"
$ws.Cells.Item(8, 4).Value = "#faae53"
$ws.Cells.Item(8, 10).Value = "pale orange"

# Row 9 - H - pale yellow
$ws.Cells.Item(9, 1).Value = "H"
$ws.Cells.Item(9, 2).Value = "#f8f391"
$ws.Cells.Item(9, 10).Value = "pale yellow"

# --- Turn the range into the "groups" table ---
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:J9"), [System.Type]::Missing, 1)
$lo.Name = "groups"
$lo.TableStyle = "TableStyleMedium2"

# --- Column J ("notes") width to roughly match the best-fit width used by the source file ---
$ws.Columns.Item(10).ColumnWidth = 11.333333333333332

# --- Give the new sheet a landscape-free (portrait) page setup like its siblings ---
$ws.PageSetup.Orientation = 1

# --- Select a sensible cell and make "Colors" the active tab (mirrors the saved workbook state) ---
$ws.Range("B8").Select()
$ws.Activate()
